$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Add a new row to the second table (Software/framework | Usage)
#    with "Google cloud storage" / "Image hosting".
# ------------------------------------------------------------------
$t = $d.Tables.Item(2)
$newRow = $t.Rows.Add()

$newRow.Cells.Item(1).Range.Text = "Google cloud storage"

# Type the new cell's text with one throw-away trailing character so
# that the later bookmark insertion point is not the very last
# position in the paragraph (inserting/bookmarking right at a
# paragraph-end boundary does not resolve reliably). We'll strip the
# placeholder character again right after placing the bookmark.
$c2 = $newRow.Cells.Item(2).Range
$c2.Text = "Image hostingX"

$c2after = $newRow.Cells.Item(2).Range
$bookmarkPos = $c2after.End - 2   # right after "...hosting", before the "X" and before the end-of-cell mark

# ------------------------------------------------------------------
# 2. Move the document's "_GoBack" bookmark from its old location
#    (after "(A0148076Y)") to the end of the new "Image hosting"
#    cell -- Bookmarks.Add on an existing name relocates it.
# ------------------------------------------------------------------
$bmRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the placeholder character now that the bookmark anchors the
# correct (safe) position; the bookmark stays put as the surrounding
# text changes.
$placeholder = $d.Range($bookmarkPos, $bookmarkPos + 1)
$placeholder.Delete()
